$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C (rows 6-9): text / formatting updates ---

# Row 6: "Pop Layer Operations" -> emphasize with red font + wrap text
$ws.Range("C6").Value = "Pop Layer Operations"
$ws.Range("C6").Font.Color = 255
$ws.Range("C6").WrapText = $true

# Row 7: "Delete" stays plain text
$ws.Range("C7").Value = "Delete"

# Row 8: re-worded string, re-entered with a quote-prefix (leading apostrophe is
# consumed by Excel as the "treat as text" marker; the literal trailing apostrophe
# in the sentence itself is preserved as part of the text)
$ws.Range("C8").Formula = "'Begin to get the list of associated subnets'"

# Row 9: re-worded string, same quote-prefix treatment
$ws.Range("C9").Formula = "'Data of the list of associated subnets are'"

# --- Sheet view: move the active selection to C10 ---
$ws.Range("C10").Select()

# --- Page setup: stamp an explicit (default) print setup ---
$ws.PageSetup.PaperSize = 0
$ws.PageSetup.Orientation = 1
